$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated "Valor_de_daño_normalizado" values (column B) resulting from the
# recalculated "DI generado" (generated damage index) function.
$values = @(
    17.262116248225382,
    23.67769743595575,
    14.68429139445821,
    12.215752612708229,
    17.727942938356442,
    3.2259170681955474,
    2.0684347964530749,
    9.4043584241022007,
    5.4872122839011297,
    9.0035282016104965,
    4.0466953965834707,
    5.563134468829988,
    4.6427438858591765,
    1.1131539804094872,
    1.2484606147144179,
    2.74993011337342,
    6.3890881919837632,
    1.9432551768959128,
    5.8742020433996984,
    2.4238759775443648,
    2.3958861441605799,
    1.0819989031569439,
    1.0381566252922283,
    1.7757158580346428,
    1.277947876667431,
    1.181466015245304,
    1.2406097640291553,
    1.4287804937813759,
    1.2983002430515234,
    1.1710920485537886,
    1.245422547318868,
    1.2458425025238651,
    1.1496706760140118,
    1.6599354397445933,
    1.3601516579964414
)

$startRow = 3
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
